# TimeSheet_Week5.xlsx update: "Team Meeting UPDATED Sep 22, 2021"
#
# Semantics of the change (reconstructed from the OOXML diff):
#   - A new task row "TA Meeting" is inserted into the weekly task table,
#     right above the existing "Organizing (misc.)" row (i.e. it becomes
#     the new row 10, pushing "Organizing (misc.)" down to row 11 and the
#     "Daily Total" row down to row 12).
#   - The hour that used to be logged for "Organizing (misc.)" on the
#     first day column (B) moves down one row with it; the new "TA
#     Meeting" row gets a fresh hour logged on the Wednesday column (D).
#   - The "Daily Total" row's column sums (H and I) widen by one row to
#     keep covering the task block (now B6:B11 .. I6:I11 instead of
#     ..B10 .. I10).
#   - A new informational line is appended at the end of the notes
#     block: "It is intended both as an accountability tool and as
#     validation for your estimates ".
#   - Cosmetic: selection moved to I13, default column width nudged from
#     8.88671875 to 8.90234375.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "TA Meeting" row above the "Organizing (misc.)" row.
# Inserting a whole row here carries formatting down with the existing
# row 10 ("Organizing (misc.)") to the new row 11, and Excel's formula
# engine widens the "Daily Total" SUM ranges (row 11 -> row 12) to keep
# including the task block, exactly like a live Insert would.
$ws.Rows("10:10").Insert()

# --- Populate the newly blank row 10 as the "TA Meeting" task.
# Match the look of the other task rows by copying their formatting in.
$ws.Range("A9:J9").Copy()
$ws.Range("A10:J10").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("J10").Style = $ws.Range("J9").Style

$ws.Range("A10").Value = "TA Meeting"
$ws.Range("I10").Formula = "=SUM(B10:H10)"

# Move the hour that used to sit in B10 ("Organizing (misc.)", now row 11)
# down with its row, and log the new TA Meeting hour on Wednesday (D10).
$ws.Range("B10").ClearContents()
$ws.Range("D10").Value = 1
$ws.Range("B11").Value = 1

# --- Make sure the Daily Total row (now row 12) sums across the full
# task block, including the freshly inserted row 11.
$ws.Range("H12").Formula = "=SUM(H6:H11)"
$ws.Range("I12").Formula = "=SUM(I6:I11)"

# Note: the informational notes block (originally rows 13-19) was already
# carried down to rows 14-20 by the row insert above, including the last
# line ("It is intended both as an accountability tool and as validation
# for your estimates "), which already existed in the workbook's shared
# strings as an as-yet-unreferenced entry - it lands on row 20 for free,
# matching the target layout with no extra write needed.

# --- Cosmetic tweaks to match the edited file: selection + default col width.
$ws.Range("I13").Select()
$ws.StandardWidth = 8.90234375
